$d = $word.ActiveDocument

$replacements = @(
    @("14×79=1106", "83×31=2573"),
    @("95×47=4465", "18×72=1296"),
    @("33×56=1848", "70×96=6720"),
    @("51×18=918", "48×68=3264"),
    @("75×24=1800", "30×62=1860"),
    @("25×99=2475", "47×65=3055"),
    @("55×57=3135", "97×48=4656"),
    @("26×61=1586", "28×55=1540"),
    @("39×97=3783", "52×73=3796"),
    @("11×20=220", "89×24=2136"),
    @("85×78=6630", "15×52=780"),
    @("54×65=3510", "52×13=676"),
    @("42×38=1596", "81×31=2511"),
    @("42×84=3528", "86×60=5160"),
    @("38×23=874", "13×19=247"),
    @("59×83=4897", "54×57=3078"),
    @("33×69=2277", "41×37=1517"),
    @("47×87=4089", "79×88=6952"),
    @("53×40=2120", "68×57=3876"),
    @("15×72=1080", "97×34=3298"),
    @("23×73=1679", "60×66=3960"),
    @("25×28=700", "41×71=2911"),
    @("32×61=1952", "17×57=969"),
    @("43×36=1548", "42×87=3654"),
    @("16×17=272", "60×44=2640")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
